$wb = $excel.ActiveWorkbook

function Set-TextValue($range, [string]$text) {
    # Force the value to be written as literal text (shared string), never
    # auto-coerced into a boolean / number / date by the engine's natural
    # "smart" value parsing. A leading apostrophe forces text entry (like
    # typing '... into Excel); resetting the style afterwards removes the
    # "quote prefix" formatting flag that the apostrophe trick leaves behind.
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

# ---------------------------------------------------------------------------
# Sheet "Overview" (sheet1.xml / table3.xml "Overview")
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A3").Value = "73367972-c8e9-41d0-a2c7-5c5ca94699f6.md"
$wsOverview.Range("B3").Value = "e2e\73367972-c8e9-41d0-a2c7-5c5ca94699f6.md"
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("G3").Value = "2016-09-07 03:04:56"

$wsOverview.Range("B3").Style = "HyperLink"
$wsOverview.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wb.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/96db7b72afa0450eb3ca973b8ba638851d7c7e01/e2e/73367972-c8e9-41d0-a2c7-5c5ca94699f6.md", "", "", "e2e\73367972-c8e9-41d0-a2c7-5c5ca94699f6.md") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "zh-cn" (sheet2.xml / table1.xml "zh-cn")
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.ListRows.Add() | Out-Null

$wsZhCn.Range("A3").Value = "73367972-c8e9-41d0-a2c7-5c5ca94699f6.md"
$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("D3").Value = "e2e"
$wsZhCn.Range("E3").Value = "ht"
Set-TextValue $wsZhCn.Range("F3") "True"
$wsZhCn.Range("G3").Value = "73367972-c8e9-41d0-a2c7-5c5ca94699f6.651483bae050d209e408959b9a71486403b6e8c7.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-09-07 03:04:44"
$wsZhCn.Range("I3").Value = "73367972-c8e9-41d0-a2c7-5c5ca94699f6.md"
$wsZhCn.Range("J3").Value = "73367972-c8e9-41d0-a2c7-5c5ca94699f6.651483bae050d209e408959b9a71486403b6e8c7.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-09-07 03:05:41"
Set-TextValue $wsZhCn.Range("L3") ""
Set-TextValue $wsZhCn.Range("M3") "True"
Set-TextValue $wsZhCn.Range("N3") ""
Set-TextValue $wsZhCn.Range("O3") "False"
Set-TextValue $wsZhCn.Range("P3") ""

$wsZhCn.Range("A3").Style = "HyperLink"
$wsZhCn.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("I3").Style = "HyperLink"
$wsZhCn.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wb.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/96db7b72afa0450eb3ca973b8ba638851d7c7e01/e2e/73367972-c8e9-41d0-a2c7-5c5ca94699f6.md", "", "", "73367972-c8e9-41d0-a2c7-5c5ca94699f6.md") | Out-Null
$wb.Hyperlinks.Add($wsZhCn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/f3952b8db7eb72545c9b582a743d2791ac23b108/e2e/73367972-c8e9-41d0-a2c7-5c5ca94699f6.md", "", "", "73367972-c8e9-41d0-a2c7-5c5ca94699f6.md") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "de-de" (sheet3.xml / table2.xml "de-de")
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.ListRows.Add() | Out-Null

$wsDeDe.Range("A3").Value = "73367972-c8e9-41d0-a2c7-5c5ca94699f6.md"
$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("D3").Value = "e2e"
$wsDeDe.Range("E3").Value = "ht"
Set-TextValue $wsDeDe.Range("F3") "True"
$wsDeDe.Range("G3").Value = "73367972-c8e9-41d0-a2c7-5c5ca94699f6.651483bae050d209e408959b9a71486403b6e8c7.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-09-07 03:04:56"
$wsDeDe.Range("I3").Value = "73367972-c8e9-41d0-a2c7-5c5ca94699f6.md"
$wsDeDe.Range("J3").Value = "73367972-c8e9-41d0-a2c7-5c5ca94699f6.651483bae050d209e408959b9a71486403b6e8c7.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-09-07 03:06:01"
Set-TextValue $wsDeDe.Range("L3") ""
Set-TextValue $wsDeDe.Range("M3") "True"
Set-TextValue $wsDeDe.Range("N3") ""
Set-TextValue $wsDeDe.Range("O3") "False"
Set-TextValue $wsDeDe.Range("P3") ""

$wsDeDe.Range("A3").Style = "HyperLink"
$wsDeDe.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("I3").Style = "HyperLink"
$wsDeDe.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wb.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/96db7b72afa0450eb3ca973b8ba638851d7c7e01/e2e/73367972-c8e9-41d0-a2c7-5c5ca94699f6.md", "", "", "73367972-c8e9-41d0-a2c7-5c5ca94699f6.md") | Out-Null
$wb.Hyperlinks.Add($wsDeDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/6dbdfb16a4ed8ab6d98955fa38cf07c0427a4ccf/e2e/73367972-c8e9-41d0-a2c7-5c5ca94699f6.md", "", "", "73367972-c8e9-41d0-a2c7-5c5ca94699f6.md") | Out-Null
